# Convert the two M2Doc field codes ("m: if self.name = 'anydsl'" and
# " m:endif ") that are stored as real Word fields
# (fldChar begin / instrText* / fldChar end) into plain text runs holding
# the literal "{m: ...}" token text -- the layout the new
# TokenIteratorFieldRewriterSplit parser expects.
#
# Each piece of the original instruction text becomes its own
# <w:r><w:t>...</w:t></w:r> run (the field delimiters are folded into the
# first/last text runs). A scratch bookmark is (re)planted at the
# insertion point after every InsertAfter() call so the host does not
# coalesce the freshly-typed, identically-formatted runs back together;
# it is removed once at the very end.

$d = $word.ActiveDocument

function Insert-Run([int]$pos, [string]$text) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $d.Bookmarks.Add("m2docSplitTmp", $d.Range($r.End, $r.End)) | Out-Null
    return $r.End
}

# --- First field: "m: if self.name = 'anydsl'" -> "{m: if self.name = 'anydsl'}"
$field1 = $d.Fields.Item(1)
$pos = $field1.Code.Start - 1
$field1.Delete()

$pieces1 = @("{m:", " ", "if ", "self.name ", "=", " ", "'", "anydsl", "'}")
foreach ($piece in $pieces1) {
    $pos = Insert-Run $pos $piece
}

# --- Second field: " m:endif " -> "{m:endif}" ---------------------------
$field2 = $d.Fields.Item(1)
$pos2 = $field2.Code.Start - 1
$field2.Delete()

$pos2 = Insert-Run $pos2 "{m:endif}"

$d.Bookmarks.Item("m2docSplitTmp").Delete()

Write-Host "Done"
